$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7464583
$ws.Range("I40").Value = 1895.5964
$ws.Range("J40").Value = 50001900
$ws.Range("K40").Value = 1895.5964
$ws.Range("L40").Value = 50001900
$ws.Range("M40").Value = -1720.5964
$ws.Range("N40").Value = -50002250
$ws.Range("H64").Value = 7214.2856
$ws.Range("I64").Value = 12166.667
$ws.Range("K64").Value = 12166.667
$ws.Range("M64").Value = -11918.667
$ws.Range("H67").Value = 7214.2856
$ws.Range("I67").Value = 12166.667
$ws.Range("K67").Value = 12166.667
$ws.Range("M67").Value = -11308.667
$ws.Range("H129").Value = 1264.64
$ws.Range("I129").Value = 815
$ws.Range("K129").Value = 2445
$ws.Range("M129").Value = 2555
$ws.Range("H132").Value = 1504
$ws.Range("I132").Value = 1067.0962
$ws.Range("J132").Value = 4749.5713
$ws.Range("K132").Value = 3201.2886
$ws.Range("L132").Value = 14248.7139
$ws.Range("M132").Value = -671.2885999999999
$ws.Range("N132").Value = -19308.7139
$ws.Range("H138").Value = 2823.611
$ws.Range("I138").Value = 1454.9259
$ws.Range("J138").Value = 4192.2964
$ws.Range("K138").Value = 4364.7777
$ws.Range("L138").Value = 12576.8892
$ws.Range("M138").Value = 775.2223000000004
$ws.Range("N138").Value = -22856.8892

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4581.39
$ws.Range("I32").Value = 4207.7837
$ws.Range("J32").Value = 16661.334
$ws.Range("K32").Value = 4207.7837
$ws.Range("L32").Value = 16661.334
$ws.Range("M32").Value = -3920.7837
$ws.Range("N32").Value = -17235.334
$ws.Range("H61").Value = 176193.23
$ws.Range("I61").Value = 4070.5264
$ws.Range("J61").Value = 503226.4
$ws.Range("K61").Value = 4070.5264
$ws.Range("L61").Value = 503226.4
$ws.Range("M61").Value = -3858.5264
$ws.Range("N61").Value = -503650.4
$ws.Range("H74").Value = 10205352
$ws.Range("I74").Value = 931.1070999999999
$ws.Range("J74").Value = 23811246
$ws.Range("K74").Value = 931.1070999999999
$ws.Range("L74").Value = 23811246
$ws.Range("M74").Value = -57.10709999999995
$ws.Range("N74").Value = -23812994
$ws.Range("H77").Value = 10205352
$ws.Range("I77").Value = 931.1070999999999
$ws.Range("J77").Value = 23811246
$ws.Range("K77").Value = 4655.5355
$ws.Range("L77").Value = 119056230
$ws.Range("M77").Value = -287.5355
$ws.Range("N77").Value = -119064966
$ws.Range("H132").Value = 1391306.6
$ws.Range("I132").Value = 1828.3954
$ws.Range("J132").Value = 3451567.2
$ws.Range("K132").Value = 5485.1862
$ws.Range("L132").Value = 10354701.6
$ws.Range("M132").Value = -2955.1862
$ws.Range("N132").Value = -10359761.6
$ws.Range("H136").Value = 176193.23
$ws.Range("I136").Value = 4070.5264
$ws.Range("J136").Value = 503226.4
$ws.Range("K136").Value = 12211.5792
$ws.Range("L136").Value = 1509679.2
$ws.Range("M136").Value = -9661.5792
$ws.Range("N136").Value = -1514779.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 19444.904
$ws.Range("I134").Value = 4097.2554
$ws.Range("J134").Value = 64528.625
$ws.Range("K134").Value = 12291.7662
$ws.Range("L134").Value = 193585.875
$ws.Range("M134").Value = -9756.7662
$ws.Range("N134").Value = -198655.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2490.5625
$ws.Range("I16").Value = 2467.7778
$ws.Range("J16").Value = 2519.8572
$ws.Range("K16").Value = 2467.7778
$ws.Range("L16").Value = 2519.8572
$ws.Range("M16").Value = -2180.7778
$ws.Range("N16").Value = -3093.8572
$ws.Range("H31").Value = 244952.8
$ws.Range("I31").Value = 1545.4286
$ws.Range("J31").Value = 538720.3
$ws.Range("K31").Value = 1545.4286
$ws.Range("L31").Value = 538720.3
$ws.Range("M31").Value = -1250.4286
$ws.Range("N31").Value = -539310.3
$ws.Range("H34").Value = 244952.8
$ws.Range("I34").Value = 1545.4286
$ws.Range("J34").Value = 538720.3
$ws.Range("K34").Value = 1545.4286
$ws.Range("L34").Value = 538720.3
$ws.Range("M34").Value = -1343.4286
$ws.Range("N34").Value = -539124.3
$ws.Range("H62").Value = 6986.3477
$ws.Range("I62").Value = 7417.5
$ws.Range("J62").Value = 6000.857
$ws.Range("K62").Value = 7417.5
$ws.Range("L62").Value = 6000.857
$ws.Range("M62").Value = -6793.5
$ws.Range("N62").Value = -7248.857
$ws.Range("H65").Value = 6986.3477
$ws.Range("I65").Value = 7417.5
$ws.Range("J65").Value = 6000.857
$ws.Range("K65").Value = 37087.5
$ws.Range("L65").Value = 30004.285
$ws.Range("M65").Value = -33967.5
$ws.Range("N65").Value = -36244.285
$ws.Range("H99").Value = 26250.5
$ws.Range("I99").Value = 26250.5
$ws.Range("K99").Value = 26250.5
$ws.Range("M99").Value = -24752.5
$ws.Range("H105").Value = 2893.5293
$ws.Range("I105").Value = 3092.1428
$ws.Range("J105").Value = 1966.6666
$ws.Range("K105").Value = 3092.1428
$ws.Range("L105").Value = 1966.6666
$ws.Range("M105").Value = -1345.1428
$ws.Range("N105").Value = -5460.6666
$ws.Range("H113").Value = 2490.5625
$ws.Range("I113").Value = 2467.7778
$ws.Range("J113").Value = 2519.8572
$ws.Range("K113").Value = 2467.7778
$ws.Range("L113").Value = 2519.8572
$ws.Range("M113").Value = -297.7777999999998
$ws.Range("N113").Value = -6859.8572
$ws.Range("H126").Value = 26250.5
$ws.Range("I126").Value = 26250.5
$ws.Range("K126").Value = 78751.5
$ws.Range("M126").Value = -76281.5
$ws.Range("H132").Value = 6669311
$ws.Range("I132").Value = 12501786
$ws.Range("J132").Value = 3625
$ws.Range("K132").Value = 37505358
$ws.Range("L132").Value = 10875
$ws.Range("M132").Value = -37502828
$ws.Range("N132").Value = -15935

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 267.11765
$ws.Range("I2").Value = 178.8
$ws.Range("J2").Value = 393.2857
$ws.Range("K2").Value = 178.8
$ws.Range("L2").Value = 393.2857
$ws.Range("M2").Value = -65.80000000000001
$ws.Range("N2").Value = -619.2857
$ws.Range("H102").Value = 2685.9092
$ws.Range("I102").Value = 2564.1304
$ws.Range("J102").Value = 2966
$ws.Range("K102").Value = 2564.1304
$ws.Range("L102").Value = 2966
$ws.Range("M102").Value = -942.1304
$ws.Range("N102").Value = -6210
$ws.Range("H126").Value = 11540
$ws.Range("I126").Value = 14320
$ws.Range("J126").Value = 3200
$ws.Range("K126").Value = 42960
$ws.Range("L126").Value = 9600
$ws.Range("M126").Value = -40490
$ws.Range("N126").Value = -14540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3572.7273
$ws.Range("I7").Value = 2285.7144
$ws.Range("J7").Value = 5825
$ws.Range("K7").Value = 2285.7144
$ws.Range("L7").Value = 5825
$ws.Range("M7").Value = -2173.7144
$ws.Range("N7").Value = -6049
$ws.Range("H40").Value = 125002610
$ws.Range("I40").Value = 200002270
$ws.Range("K40").Value = 200002270
$ws.Range("M40").Value = -200002134
$ws.Range("H68").Value = 50002012
$ws.Range("I68").Value = 1780
$ws.Range("J68").Value = 100002240
$ws.Range("K68").Value = 1780
$ws.Range("L68").Value = 100002240
$ws.Range("M68").Value = -1031
$ws.Range("N68").Value = -100003738
$ws.Range("H71").Value = 50002012
$ws.Range("I71").Value = 1780
$ws.Range("J71").Value = 100002240
$ws.Range("K71").Value = 8900
$ws.Range("L71").Value = 500011200
$ws.Range("M71").Value = -5156
$ws.Range("N71").Value = -500018688
$ws.Range("H126").Value = 3572.7273
$ws.Range("I126").Value = 2285.7144
$ws.Range("J126").Value = 5825
$ws.Range("K126").Value = 6857.1432
$ws.Range("L126").Value = 17475
$ws.Range("M126").Value = -4387.1432
$ws.Range("N126").Value = -22415
$ws.Range("H136").Value = 8845.297
$ws.Range("I136").Value = 6051.1377
$ws.Range("J136").Value = 18974.125
$ws.Range("K136").Value = 18153.4131
$ws.Range("L136").Value = 56922.375
$ws.Range("M136").Value = -15603.4131
$ws.Range("N136").Value = -62022.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1325.5
$ws.Range("I126").Value = 1170.1
$ws.Range("K126").Value = 3510.3
$ws.Range("M126").Value = -1040.3
$ws.Range("H132").Value = 1510.1628
$ws.Range("I132").Value = 791.9394
$ws.Range("J132").Value = 3880.3
$ws.Range("K132").Value = 2375.8182
$ws.Range("L132").Value = 11640.9
$ws.Range("M132").Value = 154.1818000000003
$ws.Range("N132").Value = -16700.9
$ws.Range("H136").Value = 2209.9758
$ws.Range("I136").Value = 1994.0869
$ws.Range("J136").Value = 2478.3784
$ws.Range("K136").Value = 5982.2607
$ws.Range("L136").Value = 7435.135200000001
$ws.Range("M136").Value = -3432.2607
$ws.Range("N136").Value = -12535.1352
